$d = $word.ActiveDocument

# Locate the "N° " run near the top of the document (evaluation number
# heading) and insert a new, separate run containing "00" right after it,
# matching its italic / size-24 / es-ES formatting but keeping it as a
# distinct <w:r> element (no merge with the neighbouring runs).
$rng = $d.Content
$found = $rng.Find.Execute("N° ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Collapse(0)   # wdCollapseEnd
    $rng.InsertAfter("00")

    # Force the new run to stay split from its neighbours (which would
    # otherwise be merged automatically because the formatting matches)
    # by toggling a property that is not actually part of the target
    # formatting, then reverting it.
    $rng.Font.Bold = $true
    $rng.Font.Bold = $false
}
